$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.313.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "'2.709.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.28%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'608.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").Value = "'166.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.86%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.54%  "

$ws.Range("D9").Value = "'2.708.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.28%  "

$ws.Range("D10").Value = "'0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.77%  "

$ws.Range("D12").Value = "'0.363"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.76%  "

$ws.Range("D13").Value = "'5.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").Value = "'28.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").Value = "'3.205.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.34%  "

$ws.Range("D16").Value = "'0.0000187"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").Value = "'68.219.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").Value = "'2.706.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.86%  "

$ws.Range("D19").Value = "'11.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.15%  "

$ws.Range("D20").Value = "'369.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.34%  "

$ws.Range("D21").Value = "'7.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.94%  "

$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("D23").Value = "'4.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.73%  "

$ws.Range("D24").Value = "'2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("D25").Value = "'72.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.85%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "'10.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").Value = "'576.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").Value = "'8.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").Value = "'1.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.52%  "

$ws.Range("E35").Value = "  +1.47%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "'1.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.05%  "

$ws.Range("D38").Value = "'19.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.97%  "

$ws.Range("D39").Value = "'160.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("D40").Value = "'0.378"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.04%  "

$ws.Range("E41").Value = "  +1.05%  "

$ws.Range("D42").Value = "'1.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("E43").Value = "  +0.88%  "

$ws.Range("D44").Value = "'2.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.03%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").Value = "'0.0₆0310"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.62%  "

$ws.Range("D47").Value = "'40.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.12%  "

$ws.Range("D48").Value = "'0.595"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.50%  "

$ws.Range("D49").Value = "'154.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("E50").Value = "  +1.82%  "

$ws.Range("E51").Value = "  +3.66%  "
